{"js": "// Add a new bulleted \"ListParagraph\" item after the existing\n// \"S\u1ed5 c\u00e1i h\u1ee3p \u0111\u1ed3ng mua h\u00e0ng...\" bullet, containing the note about the\n// missing exchange-rate setup between units of measure in the same group.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The target bullet is the last paragraph in the document body (the one\n// right before the closing sectPr) \u2014 grab it via getLast() so we don't\n// depend on a brittle index.\nconst lastParagraph = paragraphs.getLast();\nlastParagraph.load(\"text\");\nawait context.sync();\n\n// insertParagraph(\"\", \"After\") clones the current paragraph's pPr/rPr\n// (style \"ListParagraph\", numPr ilvl=0/numId=1, Arial 12pt) exactly like\n// pressing Enter at the end of the bullet in Word, so the new paragraph\n// naturally joins the same bulleted list.\nconst newParagraph = lastParagraph.insertParagraph(\"\", \"After\");\n\n// Insert the three runs in order, matching the authored edit.\nnewParagraph.insertText(\n  \"Thi\u1ebfu ph\u1ea7n setup t\u1ec9 l\u1ec7 quy \u0111\u1ed5i gi\u1eefa c\u00e1c \u0111\u01a1n v\u1ecb t\u00ednh c\u00f9ng \",\n  \"End\"\n);\nnewParagraph.insertText(\"nh\u00f3m\", \"End\");\nnewParagraph.insertText(\".\", \"End\");\n\nawait context.sync();\n", "ps1": "# Add a new bulleted \"ListParagraph\" item after the existing\n# \"S\u1ed5 c\u00e1i h\u1ee3p \u0111\u1ed3ng mua h\u00e0ng...\" bullet, containing the note about the\n# missing exchange-rate setup between units of measure in the same group.\n\n$d = $word.ActiveDocument\n\n# The target bullet is the last paragraph in the document body (the one\n# right before the closing sectPr).\n$lastPara = $d.Paragraphs.Last\n\n# InsertParagraphAfter() behaves like pressing Enter at the end of that\n# bullet: the new paragraph clones its pPr/rPr (style \"ListParagraph\",\n# numPr ilvl=0/numId=1, Arial 12pt), so it naturally joins the same\n# bulleted list.\n$lastPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Last\n$newRng = $newPara.Range\n\n# Insert the three runs in order, matching the authored edit.\n$newRng.InsertAfter(\"Thi\u1ebfu ph\u1ea7n setup t\u1ec9 l\u1ec7 quy \u0111\u1ed5i gi\u1eefa c\u00e1c \u0111\u01a1n v\u1ecb t\u00ednh c\u00f9ng \")\n$newRng.Collapse(0)\n$newRng.InsertAfter(\"nh\u00f3m\")\n$newRng.Collapse(0)\n$newRng.InsertAfter(\".\")\n"}
